$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the B-column values (new prediction scores)
$ws.Range("B2").Value = 4.5238964875741399
$ws.Range("B3").Value = 2.9279481433207639
$ws.Range("B4").Value = -1.4664595685736472

# Refresh formatting on the label columns/header row so they pick up a
# new (visually identical) style entry, matching the re-saved workbook.
$ws.Range("A1:C1").HorizontalAlignment = 1
$ws.Range("A2:A4").HorizontalAlignment = 1
